$wb = $excel.ActiveWorkbook

# On the "addReseller" sheet, the commLang/billLang sample values lose the
# surrounding literal quote characters: "English (US)" -> English (US)
$wsAdd = $wb.Worksheets.Item("addReseller")
$wsAdd.Range("D2").Value = "English (US)"
$wsAdd.Range("E2").Value = "English (US)"

# Move the selection on addReseller and make it the active sheet/tab
# (previously "searchReseller" was the active tab).
$wsAdd.Activate()
$wsAdd.Range("F6").Select()
